$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 641-642, shifting existing rows 641:668 down to 643:670.
$ws.Rows("641:642").Insert()

# New row 641: Alcachofa Española Primera, report date 45147 (2023-08-09)
$ws.Range("A641").Value = 9
$ws.Range("B641").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C641").Value = "Metropolitana"
$ws.Range("D641").Value = 45147
$ws.Range("E641").Value = 13
$ws.Range("F641").Value = 100112013
$ws.Range("G641").Value = "Alcachofa"
$ws.Range("H641").Value = "Española"
$ws.Range("I641").Value = "Primera"
$ws.Range("J641").Value = 52
$ws.Range("K641").Value = 14000
$ws.Range("L641").Value = 15000
$ws.Range("M641").Value = 14500
$ws.Range("N641").Value = "`$/caja 30 unidades"
$ws.Range("O641").Value = "Provincia de Limarí"
$ws.Range("P641").Value = 483
$ws.Range("Q641").Value = 30
$ws.Range("R641").Value = "Hortaliza"

# New row 642: Alcachofa Española Segunda, report date 45147 (2023-08-09)
$ws.Range("A642").Value = 9
$ws.Range("B642").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C642").Value = "Metropolitana"
$ws.Range("D642").Value = 45147
$ws.Range("E642").Value = 13
$ws.Range("F642").Value = 100112013
$ws.Range("G642").Value = "Alcachofa"
$ws.Range("H642").Value = "Española"
$ws.Range("I642").Value = "Segunda"
$ws.Range("J642").Value = 52
$ws.Range("K642").Value = 13000
$ws.Range("L642").Value = 13000
$ws.Range("M642").Value = 13000
$ws.Range("N642").Value = "`$/caja 40 unidades"
$ws.Range("O642").Value = "Provincia de Limarí"
$ws.Range("P642").Value = 325
$ws.Range("Q642").Value = 40
$ws.Range("R642").Value = "Hortaliza"
